# Adapt Szenario Default and Link
# Updates the "Power Storage" scenario data (ExisUnits / MaxInvest columns)
# and restores the last on-screen selection in the frozen (bottom-left) pane.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ExisUnits (column E)
$ws.Range("E7").Value  = 7
$ws.Range("E10").Value = 28

# MaxInvest (column S)
$ws.Range("S7").Value  = 8
$ws.Range("S8").Value  = 8
$ws.Range("S9").Value  = 8
$ws.Range("S10").Value = 8
$ws.Range("S11").Value = 8

# Move the selection in the frozen bottom-left pane to S12 (was R18)
$ws.Range("S12").Select() | Out-Null
